$wb = $excel.ActiveWorkbook

# --- "Range Variables" sheet: remove the "var0" row ---------------------
# The row holding the first range-variable entry ("var0") is deleted
# outright; Excel shifts every row below it up by one. Because each
# var-name string shifts down by exactly one slot in the shared-string
# table too (the now-unused "var0" string is dropped on save), every
# remaining row keeps displaying the *next* variable name (var1, var2, ...),
# and the very last row (which used to hold "var30") simply disappears.
$wsRange = $wb.Worksheets.Item("Range Variables")
$wsRange.Rows("2:2").Delete() | Out-Null

# Update the cursor position left on that sheet.
$wsRange.Range("D7").Select() | Out-Null

# --- "Operators" sheet: cursor moved, no longer the active tab ----------
$wsOperators = $wb.Worksheets.Item("Operators")
$wsOperators.Activate() | Out-Null
$wsOperators.Range("C3").Select() | Out-Null

# --- "Parameters" sheet: widen column A, move cursor, make it active ----
$wsParameters = $wb.Worksheets.Item("Parameters")
$wsParameters.Activate() | Out-Null
$wsParameters.Columns.Item(1).ColumnWidth = 26.6
$wsParameters.Range("B3").Select() | Out-Null
